$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8..126 down to 9..127
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new data record
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 44545
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100101
$ws.Range("H8").Value = "Berries"
$ws.Range("I8").Value = 100101004
$ws.Range("J8").Value = "Frambuesa"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 8000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("Q8").Value = '$/bandeja 2 kilos'
$ws.Range("R8").Value = "Provincia de Curicó"
$ws.Range("S8").Value = 4000
$ws.Range("T8").Value = 2
